$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Cotizacion N°" label next to the title, styled like the
# other section headers (A3 "Datos del Cliente" / A9 "Datos del Emisor").
$ws.Range("C1").Value = "Cotizacion N°"
$ws.Range("A3").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# The quotation number itself, plain/default formatting.
$ws.Range("D1").Value = 96231

# Update the current selection (previously scrolled to B22 with
# topLeftCell A13; now at D5 with the sheet scrolled back to the top).
$ws.Range("D5").Select()

$excel.CutCopyMode = $false
